$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$Bvals = @(1.369673264576136, 1.274332601997514, 1.216979766397799, 1.193902764974666, 1.190088542737669, 1.216667353257037, 1.336551240244773, 1.581241319452374, 1.76714293803775, 1.853107643595138, 1.885865978220863, 1.878801678702359, 1.855798549673693, 1.841735348139764, 1.76155343477069, 1.712725402519766, 1.684772040047221, 1.675329916733801, 1.717909616143174, 1.862549516544277, 1.958279240978698, 1.907075254358176, 1.715565464324698, 1.513990431148045)
$Cvals = @(0.2656826914772807, 0.2404017378414096, 0.2251144228771125, 0.2189427116898059, 0.2179213725983118, 0.2250309557397827, 0.2569162633849089, 0.3213673116834457, 0.369977703328459, 0.392384241964578, 0.4009127149925575, 0.3990739936657519, 0.3930850025048471, 0.3894222933198535, 0.3685194074546985, 0.3557723090223135, 0.3484680475816617, 0.3459996365424445, 0.3571263998189806, 0.394842919606333, 0.4197476772659456, 0.4064317475838379, 0.3565141400511038, 0.3037169158874349)
$Evals = @(0.1710439527779499, 0.165189469593642, 0.1616932456686939, 0.1602927282571969, 0.1600616208067507, 0.1616742604299972, 0.16900456883684, 0.1841853705781489, 0.1958657984561114, 0.2013018212315103, 0.2033785466810869, 0.2029304666403533, 0.2014723060265595, 0.2005815315018609, 0.1955130585430425, 0.1924355245067204, 0.1906769057663453, 0.1900834241915845, 0.1927619389888307, 0.2019001037901447, 0.2079789165751365, 0.2047245950608385, 0.1926143337294164, 0.1799885989145196)
$Fvals = @(4.455129512108329, 4.262079529812581, 4.14475348412256, 4.097236510153522, 4.089363861474254, 4.144111473987607, 4.388309605089347, 4.877207620964981, 5.243170703513272, 5.411273319071654, 5.47517439657571, 5.461401119688617, 5.416525541656767, 5.389070086280185, 5.232218395774282, 5.136418511398801, 5.081469173981702, 5.062890093730857, 5.146600741265388, 5.429699871057267, 5.616150400598087, 5.516503908613913, 5.141996959463228, 4.743807069104633)
$Gvals = @(0.002515153514457071, 0.002521948973448715, 0.002526332801668354, 0.002528172617966939, 0.002528481347906592, 0.002526357397657708, 0.00251745284481935, 0.002501658464652814, 0.00249105698207625, 0.002486448828903751, 0.002484734456985734, 0.002485102318479291, 0.002486307173576282, 0.002487049165694134, 0.002491362434440138, 0.002494063277134786, 0.002495636934261532, 0.002496173224124043, 0.002493773678569033, 0.002485952448767125, 0.002481019293982693, 0.002483635951275352, 0.002493904540977125, 0.00250575418646019)
$Jvals = @(0.3130241731869177, 0.2913821266990482, 0.2780906112967472, 0.2726721111909995, 0.2717722089866186, 0.2780175457129133, 0.3055618761882641, 0.359611442139169, 0.3994306975161805, 0.4175884916315908, 0.424472202089845, 0.4229893040255774, 0.4181546549871769, 0.4151943451770137, 0.3982450478927433, 0.387859448426866, 0.3818899626019174, 0.3798694467578372, 0.3889645850584884, 0.4195744866907489, 0.4396256744430218, 0.4289192994278892, 0.3884649487426657, 0.3449768447763972)
$Mvals = @(0.2883275366245215, 0.2924932427373079, 0.2955575123909995, 0.2969333593368546, 0.2971694895933261, 0.2955755530840847, 0.2896585995623688, 0.2820863622141445, 0.2789990654183399, 0.2781366988224931, 0.2778884718051202, 0.2779384417140136, 0.2781147048569039, 0.2782328841093715, 0.2790663631577104, 0.2797167773281615, 0.280141875327331, 0.2802945535915597, 0.279642258917967, 0.2780608030031999, 0.2774840057951309, 0.2777499255543958, 0.2796757893486514, 0.2837014578521178)

for ($i = 0; $i -lt 24; $i++) {
    $row = $i + 2
    $ws.Range("B" + $row).Value2 = $Bvals[$i]
    $ws.Range("C" + $row).Value2 = $Cvals[$i]
    $ws.Range("E" + $row).Value2 = $Evals[$i]
    $ws.Range("F" + $row).Value2 = $Fvals[$i]
    $ws.Range("G" + $row).Value2 = $Gvals[$i]
    $ws.Range("J" + $row).Value2 = $Jvals[$i]
    $ws.Range("M" + $row).Value2 = $Mvals[$i]
}
